# Final updates for initial submission:
#  - Replace the "Instrumental" example (D2) with the updated literature example
#  - Fill in the previously-empty "Relational" example (D4)
#  - Bump the table's font size from ~7pt to 9pt (readability pass)
#  - Resize rows to fit the newly-expanded example text
#  - Leave the selection on D5, just below the filled-in table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content updates -------------------------------------------------
$ws.Range("D2").Value = "Reduced polio transmission and reduced health cost due to improved targeting of populations (Borowitz et al. 2023); increased consumer surplus of crops due to improved weather forecasts (Cooke and Golub 2020); avoided losses from improved wildfire suppression (Herr et al. 2020)"
$ws.Range("D4").Value = "Inceased agency of Indigenous communities for monitoring and enforcing illegal deforestation (Gonzalez et al. 2023); sense of community and quality of life through common understanding of decision contexts (Sawyer et al. 2022)"

# --- formatting updates ------------------------------------------------
# Bump the whole table to 9pt (was ~7pt)
$ws.Range("A1:D4").Font.Size = 9

# Resize rows so the longer wrapped text fits
$ws.Rows(1).RowHeight = 12.75
$ws.Rows(2).RowHeight = 96.75
$ws.Rows(3).RowHeight = 60.75
$ws.Rows(4).RowHeight = 84

# --- leave selection where the author left it --------------------------
$ws.Range("D5").Select()
